# Refresh the cryptos list (Coin / Link / Price / Volume(1h)) with the
# latest scrape. Price (D) and Volume(1h) (E) columns are plain text in
# this sheet (no number format), and a handful of coins swapped rank
# order, so row B/C (name/link) occasionally changes too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Forces a value onto a cell as literal text even when it looks like a
# number (e.g. "576.69"), matching the original inline-string cells —
# Excel would otherwise silently coerce a bare "4.24" to a Number.
# ClearFormats() afterwards drops the temporary "@" override so the
# cell's style index is left exactly as it started (no stray s="...").
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$updates = @(
    @{ Row = 2;  D = '63.596.87';  E = '  +2.55%  ' },
    @{ Row = 3;  D = '2.473.39' },
    @{ Row = 4;  E = '  +0.12%  ' },
    @{ Row = 5;  D = '576.69';    E = '  +2.44%  ' },
    @{ Row = 6;  D = '148.64';    E = '  +3.75%  ' },
    @{ Row = 7;  D = '0.998';     E = '  -0.21%  ' },
    @{ Row = 8;  E = '  +1.98%  ' },
    @{ Row = 9;  E = '  +4.10%  ' },
    @{ Row = 10; E = '  +0.74%  ' },
    @{ Row = 11; E = '  +3.66%  ' },
    @{ Row = 12; E = '  +2.43%  ' },
    @{ Row = 13; D = '27.15';     E = '  +3.75%  ' },
    @{ Row = 14; E = '  +6.15%  ' },
    @{ Row = 16; D = '63.280.43'; E = '  +2.23%  ' },
    @{ Row = 17; D = '2.479.94';  E = '  +2.58%  ' },
    @{ Row = 18; D = '11.51';     E = '  +1.79%  ' },
    @{ Row = 19; E = '  +7.88%  ' },
    @{ Row = 20; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '328.93'; E = '  +1.61%  ' },
    @{ Row = 21; B = 'Polkadot';    C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot';       D = '4.24';   E = '  +2.61%  ' },
    @{ Row = 22; B = 'Dai';         C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai';             D = '1.04';   E = '  +4.35%  ' },
    @{ Row = 23; B = 'SuiNetwork';  C = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui';      D = '1.98';   E = '  +13.80%  ' },
    @{ Row = 24; D = '67.38';     E = '  +0.42%  ' },
    @{ Row = 25; D = '630.86';    E = '  +13.55%  ' },
    @{ Row = 26; D = '8.94';      E = '  +1.85%  ' },
    @{ Row = 27; E = '  +14.30%  ' },
    @{ Row = 28; D = '2.595.85';  E = '  +2.06%  ' },
    @{ Row = 29; E = '  +9.22%  ' },
    @{ Row = 30; E = '  -0.17%  ' },
    @{ Row = 31; D = '8.41';      E = '  +2.61%  ' },
    @{ Row = 33; E = '  +3.46%  ' },
    @{ Row = 34; E = '  +9.85%  ' },
    @{ Row = 35; D = '1.55';      E = '  +2.78%  ' },
    @{ Row = 36; D = '0.998';     E = '  -0.14%  ' },
    @{ Row = 37; E = '  +2.16%  ' },
    @{ Row = 38; D = '5.54';      E = '  +1.60%  ' },
    @{ Row = 39; D = '18.98';     E = '  +1.86%  ' },
    @{ Row = 40; E = '  +2.29%  ' },
    @{ Row = 41; D = '146.37';    E = '  -3.97%  ' },
    @{ Row = 42; E = '  +20.50%  ' },
    @{ Row = 43; E = '  +0.35%  ' },
    @{ Row = 44; D = '150.61';    E = '  +2.19%  ' },
    @{ Row = 45; E = '  +4.01%  ' },
    @{ Row = 46; E = '  +4.13%  ' },
    @{ Row = 47; D = '21.11';     E = '  +6.49%  ' },
    @{ Row = 48; E = '  +2.64%  ' },
    @{ Row = 49; E = '  +5.37%  ' },
    @{ Row = 50; E = '  +1.02%  ' },
    @{ Row = 51; B = 'BitgetToken'; C = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'; D = '1.08'; E = '  +0.92%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey('D')) { Set-TextValue $ws.Range("D$r") $u.D }
    if ($u.ContainsKey('E')) { $ws.Range("E$r").Value = $u.E }
}
